$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Retailer")

# Rename existing headers
$ws.Range("A1").Value = "retailerId"
$ws.Range("B1").Value = "retailerUsername"
$ws.Range("C1").Value = "retailerPwd"
$ws.Range("D1").Value = "sfscCaseBrandName"
$ws.Range("H1").Value = "EcomLocations"
$ws.Range("J1").Value = "brandName"

# New headers for gift card columns
$ws.Range("O1").Value = "egcSku"
$ws.Range("P1").Value = "egcLocations"
$ws.Range("Q1").Value = "egcCtryCode"
$ws.Range("R1").Value = "egcPgmGroupName"
$ws.Range("S1").Value = "epgcSku"
$ws.Range("T1").Value = "epgcLocations"

# New data row 2 (AE)
$ws.Range("O2").Value = 100000002005
$ws.Range("O2").NumberFormat = "0"
$ws.Range("P2").Value = "EGIFT_107"
$ws.Range("Q2").Value = "ae"
$ws.Range("R2").Value = "TUMI UAE E-Gift Card"
$ws.Range("S2").Value = 100000002007
$ws.Range("S2").NumberFormat = "0"
$ws.Range("T2").Value = "EPGC_107"

# New data row 3 (SA)
$ws.Range("Q3").Value = "sa"
$ws.Range("R3").Value = "TUMI KSA E-Gift Card"

# Make Retailer the active sheet/tab and select R1 like the target workbook
$ws.Activate()
$ws.Range("R1").Select() | Out-Null
